$d = $word.ActiveDocument

$pairs = @(
    @("88×93=", "23×13="),
    @("42×57=", "92×92="),
    @("32×93=", "49×15="),
    @("13×52=", "36×30="),
    @("37×37=", "46×47="),
    @("73×27=", "81×42="),
    @("22×58=", "91×86="),
    @("83×81=", "39×48="),
    @("90×85=", "16×15="),
    @("40×15=", "92×99="),
    @("26×87=", "32×82="),
    @("43×16=", "54×50="),
    @("41×19=", "37×80="),
    @("28×91=", "67×20="),
    @("52×95=", "58×88="),
    @("30×11=", "24×86="),
    @("93×95=", "14×67="),
    @("29×51=", "80×12="),
    @("60×68=", "91×31="),
    @("30×54=", "65×82="),
    @("36×35=", "32×53="),
    @("24×51=", "20×33="),
    @("82×46=", "27×88="),
    @("23×27=", "72×53="),
    @("14×36=", "65×90=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
